# Replace decimal point with decimal comma in percentage values
# (only the entries that contain a "." decimal separator get changed;
#  plain integer percentages like 60%, 49%, 52%, 100% stay as-is)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "53,1%"
$ws.Range("E2").Value = "46,9%"
$ws.Range("G2").Value = "13,3%"

$ws.Range("C3").Value = "57,7%"
$ws.Range("E3").Value = "42,3%"
$ws.Range("G3").Value = "1,5%"

$ws.Range("C4").Value = "49,1%"
$ws.Range("E4").Value = "50,9%"
$ws.Range("G4").Value = "36,9%"

$ws.Range("C5").Value = "37,9%"
$ws.Range("E5").Value = "62,1%"
$ws.Range("G5").Value = "0,5%"

$ws.Range("G7").Value = "38,4%"

$ws.Range("G8").Value = "9,4%"
$ws.Range("C9").Value = "50,7%"
$ws.Range("E9").Value = "49,3%"
